$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new observation was recorded for 2026/02/27 (Friday) at time-slot 8,
# ranked 33. It belongs right after the existing 2026/02/27 rows (867)
# and before the 2026/12/29 block (previously row 868), so insert a new
# row at 868 and shift everything below it down by one.
$ws.Rows.Item(868).Insert()

# Column A stores the date as literal text (e.g. "2026/02/27"), not a
# real date value. Force the cell to text format first so Excel's
# automatic type detection doesn't turn the string into a date serial,
# then clear the temporary formatting so no stray style is left behind.
$ws.Range("A868").NumberFormat = "@"
$ws.Range("A868").Value = "2026/02/27"
$ws.Range("A868").ClearFormats()

$ws.Range("B868").Value = "金"
$ws.Range("C868").Value = 8
$ws.Range("D868").Value = 33
